# Load the population-weighted summary table by group (on-track, off-track)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (data starts at row 2, header at row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Map the legacy status labels to the new on-track / off-track categories
$statusMap = @{
    "Achieved"             = "on-track"
    "On Track"             = "on-track"
    "Acceleration Needed"  = "off-track"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $current = $ws.Cells.Item($r, 2).Value2
    if ($null -ne $current -and $statusMap.ContainsKey($current)) {
        $ws.Cells.Item($r, 2).Value2 = $statusMap[$current]
    }

    # Kosovo's ISO code changed from the temporary "RKS" to "XKX"
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -eq "RKS") {
        $ws.Cells.Item($r, 1).Value2 = "XKX"
    }
}
